# Apply scheduled market-price refresh updates to each profession sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 697.0526
$ws.Range("J17").Value = 720.2353000000001
$ws.Range("L17").Value = 2160.7059
$ws.Range("N17").Value = -2496.7059
# Row 29
$ws.Range("H29").Value = 3442.6667
$ws.Range("J29").Value = 4572
$ws.Range("L29").Value = 13716
$ws.Range("N29").Value = -14278
# Row 61
$ws.Range("H61").Value = 268.33334
$ws.Range("I61").Value = 268.33334
$ws.Range("K61").Value = 805.0000200000001
$ws.Range("M61").Value = -633.0000200000001
# Row 98
$ws.Range("H98").Value = 1218.8
$ws.Range("I98").Value = 700
$ws.Range("J98").Value = 1997
$ws.Range("K98").Value = 700
$ws.Range("L98").Value = 1997
$ws.Range("M98").Value = 798
$ws.Range("N98").Value = -4993
# Row 122
$ws.Range("H122").Value = 1218.8
$ws.Range("I122").Value = 700
$ws.Range("J122").Value = 1997
$ws.Range("K122").Value = 2100
$ws.Range("L122").Value = 5991
$ws.Range("M122").Value = 350
$ws.Range("N122").Value = -10891
# Row 125
$ws.Range("H125").Value = 1599.8334
$ws.Range("I125").Value = 1599.75
$ws.Range("J125").Value = 1600
$ws.Range("K125").Value = 14397.75
$ws.Range("L125").Value = 14400
$ws.Range("M125").Value = -11937.75
$ws.Range("N125").Value = -19320

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 180.76923
$ws.Range("I4").Value = 192.95653
$ws.Range("J4").Value = 87.333336
$ws.Range("K4").Value = 192.95653
$ws.Range("L4").Value = 87.333336
$ws.Range("M4").Value = -76.95652999999999
$ws.Range("N4").Value = -319.333336
# Row 32
$ws.Range("H32").Value = 2634921.5
$ws.Range("I32").Value = 3093.2222
$ws.Range("K32").Value = 3093.2222
$ws.Range("M32").Value = -2806.2222
# Row 44
$ws.Range("H44").Value = 12972.75
$ws.Range("J44").Value = 12972.75
$ws.Range("L44").Value = 12972.75
$ws.Range("N44").Value = -13948.75
# Row 122
$ws.Range("H122").Value = 2089.5
$ws.Range("I122").Value = 1610.8
$ws.Range("K122").Value = 4832.4
$ws.Range("M122").Value = -2382.4

$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 1049
$ws.Range("I64").Value = 1023.5
$ws.Range("J64").Value = 1100
$ws.Range("K64").Value = 1023.5
$ws.Range("L64").Value = 1100
$ws.Range("M64").Value = -798.5
$ws.Range("N64").Value = -1550
# Row 67
$ws.Range("H67").Value = 1049
$ws.Range("I67").Value = 1023.5
$ws.Range("J67").Value = 1100
$ws.Range("K67").Value = 1023.5
$ws.Range("L67").Value = 1100
$ws.Range("M67").Value = -243.5
$ws.Range("N67").Value = -2660

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 3501.3333
$ws.Range("I7").Value = 5474
$ws.Range("J7").Value = 94
$ws.Range("K7").Value = 5474
$ws.Range("L7").Value = 94
$ws.Range("M7").Value = -5361
$ws.Range("N7").Value = -320
# Row 12
$ws.Range("H12").Value = 15
$ws.Range("I12").Value = 15
$ws.Range("K12").Value = 15
$ws.Range("M12").Value = 155
# Row 59
$ws.Range("H59").Value = 197637620
$ws.Range("J59").Value = 197637620
$ws.Range("L59").Value = 197637620
$ws.Range("N59").Value = -197639910
# Row 99
$ws.Range("H99").Value = 2283.625
$ws.Range("I99").Value = 2257
$ws.Range("J99").Value = 2328
$ws.Range("K99").Value = 2257
$ws.Range("L99").Value = 2328
$ws.Range("M99").Value = -759
$ws.Range("N99").Value = -5324
# Row 122
$ws.Range("H122").Value = 1705.3334
$ws.Range("I122").Value = 1055.5
$ws.Range("J122").Value = 3005
$ws.Range("K122").Value = 3166.5
$ws.Range("L122").Value = 9015
$ws.Range("M122").Value = -716.5
$ws.Range("N122").Value = -13915
# Row 126
$ws.Range("H126").Value = 2283.625
$ws.Range("I126").Value = 2257
$ws.Range("J126").Value = 2328
$ws.Range("K126").Value = 6771
$ws.Range("L126").Value = 6984
$ws.Range("M126").Value = -4301
$ws.Range("N126").Value = -11924
# Row 132
$ws.Range("H132").Value = 1599.6774
$ws.Range("I132").Value = 1542.1428
$ws.Range("J132").Value = 2136.6667
$ws.Range("K132").Value = 4626.428400000001
$ws.Range("L132").Value = 6410.000100000001
$ws.Range("M132").Value = -2096.428400000001
$ws.Range("N132").Value = -11470.0001

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 992
$ws.Range("I3").Value = 992
$ws.Range("K3").Value = 2976
$ws.Range("M3").Value = -2864
# Row 92
$ws.Range("H92").Value = 340
$ws.Range("J92").Value = 200
$ws.Range("L92").Value = 600
$ws.Range("N92").Value = -3096
# Row 103
$ws.Range("H103").Value = 472.125
$ws.Range("I103").Value = 311.75
$ws.Range("J103").Value = 632.5
$ws.Range("K103").Value = 935.25
$ws.Range("L103").Value = 1897.5
$ws.Range("M103").Value = -56.25
$ws.Range("N103").Value = -3655.5
# Row 132
$ws.Range("H132").Value = 2290.818
$ws.Range("I132").Value = 1957.1428
$ws.Range("J132").Value = 2874.75
$ws.Range("K132").Value = 17614.2852
$ws.Range("L132").Value = 25872.75
$ws.Range("M132").Value = -15084.2852
$ws.Range("N132").Value = -30932.75

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 91.8
$ws.Range("I2").Value = 52.81818
$ws.Range("K2").Value = 52.81818
$ws.Range("M2").Value = 60.18182
# Row 80
$ws.Range("H80").Value = 10184.286
$ws.Range("I80").Value = 8259.200000000001
$ws.Range("K80").Value = 8259.200000000001
$ws.Range("M80").Value = -7261.200000000001
# Row 83
$ws.Range("H83").Value = 10184.286
$ws.Range("I83").Value = 8259.200000000001
$ws.Range("K83").Value = 41296
$ws.Range("M83").Value = -36304
# Row 102
$ws.Range("H102").Value = 1385
$ws.Range("I102").Value = 1385
$ws.Range("K102").Value = 1385
$ws.Range("M102").Value = 237
# Row 113
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -2830
# Row 122
$ws.Range("H122").Value = 1260.3334
$ws.Range("I122").Value = 1357.7142
$ws.Range("J122").Value = 919.5
$ws.Range("K122").Value = 4073.1426
$ws.Range("L122").Value = 2758.5
$ws.Range("M122").Value = -1623.1426
$ws.Range("N122").Value = -7658.5
# Row 123
$ws.Range("H123").Value = 975000
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
# Row 126
$ws.Range("H126").Value = 2839.8
$ws.Range("I126").Value = 2999.6667
$ws.Range("K126").Value = 8999.000100000001
$ws.Range("M126").Value = -6529.000100000001
# Row 132
$ws.Range("H132").Value = 3651
$ws.Range("I132").Value = 3299.7
$ws.Range("J132").Value = 4236.5
$ws.Range("K132").Value = 9899.099999999999
$ws.Range("L132").Value = 12709.5
$ws.Range("M132").Value = -7369.099999999999
$ws.Range("N132").Value = -17769.5

$ws = $wb.Worksheets.Item("LTW")
# Row 20
$ws.Range("H20").Value = 507499.5
$ws.Range("I20").Value = 15000
$ws.Range("K20").Value = 15000
$ws.Range("M20").Value = -14774
# Row 22
$ws.Range("H22").Value = 945.3333
$ws.Range("I22").Value = 934.4
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 934.4
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -639.4
$ws.Range("N22").Value = -1590
# Row 27
$ws.Range("H27").Value = 945.3333
$ws.Range("I27").Value = 934.4
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 934.4
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -827.4
$ws.Range("N27").Value = -1214
# Row 40
$ws.Range("H40").Value = 2000
$ws.Range("I40").Value = 2000
$ws.Range("K40").Value = 2000
$ws.Range("M40").Value = -1864
# Row 46
$ws.Range("H46").Value = 5376
$ws.Range("I46").Value = 3180
$ws.Range("K46").Value = 3180
$ws.Range("M46").Value = -2992
# Row 55
$ws.Range("H55").Value = 1572.7059
$ws.Range("I55").Value = 1026.3334
$ws.Range("J55").Value = 2187.375
$ws.Range("K55").Value = 1026.3334
$ws.Range("L55").Value = 2187.375
$ws.Range("M55").Value = -853.3334
$ws.Range("N55").Value = -2533.375
# Row 122
$ws.Range("H122").Value = 2328.2856
$ws.Range("I122").Value = 2059.8
$ws.Range("K122").Value = 6179.400000000001
$ws.Range("M122").Value = -3729.400000000001
# Row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
# Row 136
$ws.Range("H136").Value = 3000.111
$ws.Range("I136").Value = 2250.125
$ws.Range("K136").Value = 6750.375
$ws.Range("M136").Value = -4200.375

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1252
$ws.Range("I96").Value = 1252.25
$ws.Range("J96").Value = 1250
$ws.Range("K96").Value = 1252.25
$ws.Range("L96").Value = 1250
$ws.Range("M96").Value = 120.75
$ws.Range("N96").Value = -3996
# Row 126
$ws.Range("H126").Value = 4508.5835
$ws.Range("I126").Value = 2831.6875
$ws.Range("K126").Value = 8495.0625
$ws.Range("M126").Value = -6025.0625

